$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1090, pushing the existing 1090:1187 block down to 1092:1189
$ws.Range("A1090:A1091").EntireRow.Insert()

# New row 1090: weekly update - Cebolla, Morada(o), 1a (cosecha), Arica y Parinacota
$ws.Cells.Item(1090, 1).Value = 8
$ws.Cells.Item(1090, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1090, 3).Value = "Coquimbo"
$ws.Cells.Item(1090, 4).Value = 45166
$ws.Cells.Item(1090, 5).Value = 4
$ws.Cells.Item(1090, 6).Value = 100112004
$ws.Cells.Item(1090, 7).Value = "Cebolla"
$ws.Cells.Item(1090, 8).Value = "Morada(o)"
$ws.Cells.Item(1090, 9).Value = "1a (cosecha)"
$ws.Cells.Item(1090, 10).Value = 2000
$ws.Cells.Item(1090, 11).Value = 13500
$ws.Cells.Item(1090, 12).Value = 14000
$ws.Cells.Item(1090, 13).Value = 13750
$ws.Cells.Item(1090, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1090, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1090, 16).Value = 764
$ws.Cells.Item(1090, 17).Value = 18
$ws.Cells.Item(1090, 18).Value = "Hortaliza"

# New row 1091: weekly update - Cebolla, Sin especificar, 1a (guarda), O'Higgins
$ws.Cells.Item(1091, 1).Value = 8
$ws.Cells.Item(1091, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1091, 3).Value = "Coquimbo"
$ws.Cells.Item(1091, 4).Value = 45166
$ws.Cells.Item(1091, 5).Value = 4
$ws.Cells.Item(1091, 6).Value = 100112004
$ws.Cells.Item(1091, 7).Value = "Cebolla"
$ws.Cells.Item(1091, 8).Value = "Sin especificar"
$ws.Cells.Item(1091, 9).Value = "1a (guarda)"
$ws.Cells.Item(1091, 10).Value = 2000
$ws.Cells.Item(1091, 11).Value = 12000
$ws.Cells.Item(1091, 12).Value = 13000
$ws.Cells.Item(1091, 13).Value = 12500
$ws.Cells.Item(1091, 14).Value = "`$/malla 16 kilos"
$ws.Cells.Item(1091, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1091, 16).Value = 781
$ws.Cells.Item(1091, 17).Value = 16
$ws.Cells.Item(1091, 18).Value = "Hortaliza"
